# Ring.pptx — "Use https in the website URL."
#
# Slide 1 (the title slide) shows the project homepage as
#   http://ring-lang.net
# Update it to use the secure scheme:
#   https://ring-lang.net
#
# We locate the run by its current text and replace only the characters
# that make up the URL, leaving every other run / line-break / formatting
# attribute in the title text box untouched.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$oldUrl = "http://ring-lang.net"
$newUrl = "https://ring-lang.net"

$fullText = $tr.Text
$startPos = $fullText.IndexOf($oldUrl)

if ($startPos -ge 0) {
    $urlRange = $tr.Characters($startPos + 1, $oldUrl.Length)
    $urlRange.Text = $newUrl
}
